$d = $word.ActiveDocument

# The first paragraph currently holds a single run:
#   "This is a Microsoft word document."
# We need to turn that into four separate runs:
#   "This is a Microsoft word document." / " (" / "Changed main" / ")"
# so the paragraph text reads:
#   "This is a Microsoft word document. (Changed main)"

$firstPara = $d.Paragraphs(1).Range

# Range covering just the run's text, not the paragraph mark, so the
# paragraph's own properties (pPr, paraId, rsids, ...) are left untouched.
$runRange = $d.Range($firstPara.Start, $firstPara.End - 1)

$packageXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t>This is a Microsoft word document.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> (</w:t></w:r>
            <w:r><w:t>Changed main</w:t></w:r>
            <w:r><w:t>)</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$runRange.InsertXML($packageXml)
